# Apply cryptos list update (diff dated Mon Dec 25 09:22:08 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text content looks like a plain number (e.g. "0.999", "9.16").
# Force them to Text format first so Excel keeps them as strings (matching the
# original inlineStr/text cell type) instead of auto-converting to numeric values.
$textCells = @("D5","D6","D9","D10","D11","D12","D16","D20","D22","D23","D24","D25","D27","D29","D30","D32","D33","D34","D35","D36","D39","D41","D42","D43","D44","D46","D48","D49","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "43.176.46"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").Value = "2.276.09"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "111.44"
$ws.Range("E5").Value = "  -2.85%  "
$ws.Range("D6").Value = "264.01"
$ws.Range("E6").Value = "  -1.74%  "
$ws.Range("E7").Value = "  +2.74%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.606"
$ws.Range("E9").Value = "  -2.20%  "
$ws.Range("D10").Value = "46.91"
$ws.Range("E10").Value = "  -2.86%  "
$ws.Range("D11").Value = "0.0934"
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").Value = "9.16"
$ws.Range("E12").Value = "  +4.41%  "
$ws.Range("E13").Value = "  +1.57%  "
$ws.Range("E14").Value = "  -1.60%  "
$ws.Range("D15").Value = "2.617.53"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Value = "0.863"
$ws.Range("E16").Value = "  +2.09%  "
$ws.Range("D17").Value = "2.269.23"
$ws.Range("D18").Value = "43.149.92"
$ws.Range("E18").Value = "  -1.19%  "
$ws.Range("D20").Value = "6.76"
$ws.Range("E20").Value = "  +3.89%  "
$ws.Range("E21").Value = "  -0.92%  "
$ws.Range("D22").Value = "2.44"
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("D23").Value = "234.17"
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("D24").Value = "9.43"
$ws.Range("E24").Value = "  -3.91%  "
$ws.Range("D25").Value = "2.85"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("E26").Value = "  +2.00%  "
$ws.Range("D27").Value = "11.33"
$ws.Range("E27").Value = "  -2.60%  "
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").Value = "40.77"
$ws.Range("E29").Value = "  -1.86%  "
$ws.Range("D30").Value = "3.34"
$ws.Range("E30").Value = "  -1.69%  "
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("D32").Value = "173.05"
$ws.Range("E32").Value = "  -2.18%  "
$ws.Range("D33").Value = "21.48"
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("D34").Value = "0.0900"
$ws.Range("E34").Value = "  -3.06%  "
$ws.Range("D35").Value = "5.63"
$ws.Range("E35").Value = "  +1.18%  "
$ws.Range("D36").Value = "0.129"
$ws.Range("E36").Value = "  +1.82%  "
$ws.Range("E37").Value = "  +3.75%  "
$ws.Range("E38").Value = "  -2.46%  "
$ws.Range("D39").Value = "3.96"
$ws.Range("E39").Value = "  +3.95%  "
$ws.Range("E40").Value = "  -4.78%  "
$ws.Range("D41").Value = "2.60"
$ws.Range("E41").Value = "  +7.11%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "14.27"
$ws.Range("E42").Value = "  +2.22%  "
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").Value = "76.35"
$ws.Range("E43").Value = "  +5.63%  "
$ws.Range("D44").Value = "0.237"
$ws.Range("E44").Value = "  -2.53%  "
$ws.Range("E45").Value = "  -1.81%  "
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("E47").Value = "  -4.24%  "
$ws.Range("D48").Value = "8.53"
$ws.Range("E48").Value = "  -2.31%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").Value = "1.26"
$ws.Range("E49").Value = "  +1.97%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "101.72"
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("E51").Value = "  -0.88%  "
